$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: text row ("a","b","y","y","y","y")
$ws.Range("A4").Value = "a"
$ws.Range("B4").Value = "b"
$ws.Range("C4").Value = "y"
$ws.Range("D4").Value = "y"
$ws.Range("E4").Value = "y"
$ws.Range("F4").Value = "y"

# Row 5: numeric fraction in A5, then "b","y","y","y","y"
$ws.Range("A5").Value = 0.33333333333333331
$ws.Range("A5").NumberFormat = "???/???"
$ws.Range("B5").Value = "b"
$ws.Range("C5").Value = "y"
$ws.Range("D5").Value = "y"
$ws.Range("E5").Value = "y"
$ws.Range("F5").Value = "y"
